$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 2 and 3, shifting existing row 2 down to row 4
$ws.Range("2:3").Insert()

# Clear inherited formatting on the newly inserted rows so cells have no explicit style
$ws.Range("A2:BB3").Style = "Normal"

# Row 2 data (Id f9UtWQkR)
$ws.Range("A2").Value = "f9UtWQkR"
$ws.Range("B2").Value = "27/11/2024"
$ws.Range("C2").Value = "11:00"
$ws.Range("D2").Value = "GEORGIA - CRYSTALBET EROVNULI LIGA"
$ws.Range("E2").Value = "Telavi"
$ws.Range("F2").Value = "Samgurali"
$ws.Range("G2").Value = 2.52
$ws.Range("H2").Value = 2.92
$ws.Range("I2").Value = 2.7
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 2.02
$ws.Range("L2").Value = 3.3
$ws.Range("M2").Value = 1.04
$ws.Range("N2").Value = 6.45
$ws.Range("O2").Value = 1.29
$ws.Range("P2").Value = 2.92
$ws.Range("Q2").Value = 1.9
$ws.Range("R2").Value = 1.72
$ws.Range("S2").Value = 1.44
$ws.Range("T2").Value = 2.3
$ws.Range("U2").Value = 1.74
$ws.Range("V2").Value = 2.04
$ws.Range("W2").Value = 7.4
$ws.Range("X2").Value = 11.25
$ws.Range("Y2").Value = 7.8
$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 16.5
$ws.Range("AB2").Value = 21
$ws.Range("AC2").Value = 8.75
$ws.Range("AD2").Value = 5.1
$ws.Range("AE2").Value = 10.25
$ws.Range("AF2").Value = 40
$ws.Range("AG2").Value = 250
$ws.Range("AH2").Value = 7
$ws.Range("AI2").Value = 11.25
$ws.Range("AJ2").Value = 8.25
$ws.Range("AK2").Value = 26
$ws.Range("AL2").Value = 19
$ws.Range("AM2").Value = 24
$ws.Range("AN2").Value = 4.55
$ws.Range("AO2").Value = 13
$ws.Range("AP2").Value = 18
$ws.Range("AQ2").Value = 55
$ws.Range("AR2").Value = 75
$ws.Range("AS2").Value = 200
$ws.Range("AT2").Value = 2.55
$ws.Range("AU2").Value = 6.3
$ws.Range("AV2").Value = 50
$ws.Range("AW2").Value = 4.75
$ws.Range("AX2").Value = 15
$ws.Range("AY2").Value = 21
$ws.Range("AZ2").Value = 70
$ws.Range("BA2").Value = 100
$ws.Range("BB2").Value = 200

# Row 3 data (Id niBRS6GM)
$ws.Range("A3").Value = "niBRS6GM"
$ws.Range("B3").Value = "27/11/2024"
$ws.Range("C3").Value = "11:00"
$ws.Range("D3").Value = "INDIA - ISL"
$ws.Range("E3").Value = "Mohammedan"
$ws.Range("F3").Value = "Bengaluru FC"
$ws.Range("G3").Value = 3.7
$ws.Range("H3").Value = 3.5
$ws.Range("I3").Value = 1.95
$ws.Range("J3").Value = 4.33
$ws.Range("K3").Value = 2.2
$ws.Range("L3").Value = 2.6
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 11
$ws.Range("O3").Value = 1.29
$ws.Range("P3").Value = 3.5
$ws.Range("Q3").Value = 1.9
$ws.Range("R3").Value = 1.9
$ws.Range("S3").Value = 1.4
$ws.Range("T3").Value = 2.75
$ws.Range("U3").Value = 1.73
$ws.Range("V3").Value = 2
$ws.Range("W3").Value = 11
$ws.Range("X3").Value = 19
$ws.Range("Y3").Value = 13
$ws.Range("Z3").Value = 41
$ws.Range("AA3").Value = 29
$ws.Range("AB3").Value = 34
$ws.Range("AC3").Value = 11
$ws.Range("AD3").Value = 6.5
$ws.Range("AE3").Value = 15
$ws.Range("AF3").Value = 51
$ws.Range("AG3").Value = 201
$ws.Range("AH3").Value = 7.5
$ws.Range("AI3").Value = 9.5
$ws.Range("AJ3").Value = 9
$ws.Range("AK3").Value = 17
$ws.Range("AL3").Value = 15
$ws.Range("AM3").Value = 26
$ws.Range("AN3").Value = 5.5
$ws.Range("AO3").Value = 21
$ws.Range("AP3").Value = 29
$ws.Range("AQ3").Value = 67
$ws.Range("AR3").Value = 81
$ws.Range("AS3").Value = 201
$ws.Range("AT3").Value = 2.75
$ws.Range("AU3").Value = 8
$ws.Range("AV3").Value = 51
$ws.Range("AW3").Value = 4
$ws.Range("AX3").Value = 10
$ws.Range("AY3").Value = 21
$ws.Range("AZ3").Value = 34
$ws.Range("BA3").Value = 51
$ws.Range("BB3").Value = 151
